$wb = $excel.ActiveWorkbook

# --- Par_InitialCapacityInstalled ---
$ws = $wb.Worksheets.Item("Par_InitialCapacityInstalled")
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0

# --- Par_MaxChargeDischargeRate ---
$ws = $wb.Worksheets.Item("Par_MaxChargeDischargeRate")
$ws.Range("B7").Value = 0

# --- Par_MaxStorageCapacity ---
$ws = $wb.Worksheets.Item("Par_MaxStorageCapacity")
$ws.Range("B7").Value = 0

# --- Par_MaxDwnShift: B4:B27 all become 0.2 ---
$ws = $wb.Worksheets.Item("Par_MaxDwnShift")
for ($r = 4; $r -le 27; $r++) {
    $ws.Range("B$r").Value = 0.2
}

# --- Par_MaxUpShift: B4:B27 all become 0.05 ---
$ws = $wb.Worksheets.Item("Par_MaxUpShift")
for ($r = 4; $r -le 27; $r++) {
    $ws.Range("B$r").Value = 0.05
}

# --- Selections / active-cell updates on each touched sheet ---
$ws = $wb.Worksheets.Item("Par_InitialCapacityInstalled")
$ws.Activate()
$ws.Range("D7").Select()

$ws = $wb.Worksheets.Item("Par_MaxChargeDischargeRate")
$ws.Activate()
$ws.Range("D17").Select()

$ws = $wb.Worksheets.Item("Par_MaxStorageCapacity")
$ws.Activate()
$ws.Range("C9").Select()

$ws = $wb.Worksheets.Item("Par_MaxDwnShift")
$ws.Activate()
$ws.Range("B4:B27").Select()

$ws = $wb.Worksheets.Item("Par_MaxUpShift")
$ws.Activate()
$ws.Range("B4:B27").Select()

# --- Set_of_EnergyCarrier becomes the active/selected tab ---
$ws = $wb.Worksheets.Item("Set_of_EnergyCarrier")
$ws.Activate()
$ws.Range("P29").Select()
